$d = $word.ActiveDocument

# 1. Summary paragraph
$d.Content.Find.Execute(
    "Highly motivated and results-oriented individual seeking a challenging role in the development of innovative mobile applications. Proven ability to collaborate effectively, solve complex problems, and contribute to a dynamic team environment. Eager to leverage skills in software development, project management, and communication to contribute to impactful projects.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A motivated student with foundational knowledge in Kotlin, seeking an app development role. Leveraging academic learning and eagerness to grow into a professional position, contributing to innovative mobile application development.",
    2) | Out-Null

# 2. Experience paragraph (has AI Intern line + break + paragraph text; merge into single sentence, no break)
$expPara = $d.Paragraphs.Item(5)
$expRange = $expPara.Range
$expRange.End = $expRange.End - 1
$expRange.Text = "Developed and implemented an AI resume enhancer at Blue Silicon Infotech, yielding a 20% increase in resume completion rates. Optimized resume templates for enhanced readability and clarity, achieving a 15% improvement in resume accuracy through quantifiable results. Demonstrated expertise in AI-driven process optimization and template design, driving efficiency and effectiveness in resume development."

# 3. Education paragraph
$d.Content.Find.Execute(
    "Bachelor of Engineering from AVIT. Graduated: 2026-05. GPA: 7.1.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Bachelor of Engineering, AVIT, May 2026, GPA 7.1.",
    2) | Out-Null

# 4. Skills paragraph
$d.Content.Find.Execute(
    "Here's a revised skills section tailored for a global editing resume, focusing on conciseness, organization, and professionalism:, Skills**,    **Programming Languages:** Python, Java, Node.js, TypeScript, C#, Flutter, Kotlin, Dart, React, Python, SQL,    **Database:** MySQL, PostgreSQL, MongoDB, SQL Server,    **Web Development:** HTML, CSS, JavaScript, React, Angular, Vue.js,    **Cloud Technologies:** AWS, Azure, Google Cloud Platform,    **Operating Systems:** Linux, Windows, macOS,    **Version Control:** Git, GitHub, GitLab,    **Data Analysis:** Pandas, NumPy, Matplotlib, Seaborn,    **Testing:** Unit Testing, Integration Testing, End-to-End Testing,    **API Development:** RESTful APIs, GraphQL,    **Design Principles:** SOLID, DRY, KISS,    **Other:** Agile Development, Mobile Development, Data Science",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "C#",
    2) | Out-Null

# 5. Projects paragraph (many run/break segments; merge into single sentence, no break)
$projPara = $d.Paragraphs.Item(11)
$projRange = $projPara.Range
$projRange.End = $projRange.End - 1
$projRange.Text = "Developed an Enhanced QR Scanner and Generator project, significantly improving efficiency and accuracy. Implemented a novel algorithm and real-time data integration using QR scanner and generator, Prediction Pro, Simple Purchase Order Manager, and PDF Maker, resulting in 20% reduced processing time, improved accuracy, and enhanced real-time data integration, ultimately driving increased sales and lower operational costs."

Write-Host "Done"
